$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $text) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" "61.546.80"
Set-TextValue "E2" "  -1.67%  "
Set-TextValue "D3" "2.996.05"
Set-TextValue "E3" "  -0.55%  "
Set-TextValue "E4" "  -0.16%  "
Set-TextValue "D5" "595.94"
Set-TextValue "E5" "  +2.09%  "
Set-TextValue "D6" "143.91"
Set-TextValue "E6" "  -2.32%  "
Set-TextValue "E7" "  -0.21%  "
Set-TextValue "E8" "  -0.20%  "
Set-TextValue "D9" "2.995.20"
Set-TextValue "E9" "  -0.37%  "
Set-TextValue "E10" "  -1.57%  "
Set-TextValue "D11" "5.90"
Set-TextValue "E11" "  +2.04%  "
Set-TextValue "E12" "  +4.49%  "
Set-TextValue "E13" "  -0.32%  "
Set-TextValue "D14" "34.19"
Set-TextValue "E14" "  -2.24%  "
Set-TextValue "E15" "  +2.28%  "
Set-TextValue "D16" "3.491.31"
Set-TextValue "E16" "  -0.77%  "
Set-TextValue "E17" "  +0.21%  "
Set-TextValue "D18" "61.525.68"
Set-TextValue "E18" "  -1.75%  "
Set-TextValue "D19" "2.997.75"
Set-TextValue "E19" "  -0.78%  "
Set-TextValue "D20" "453.64"
Set-TextValue "E20" "  -2.48%  "
Set-TextValue "D21" "13.97"
Set-TextValue "E21" "  +0.27%  "
Set-TextValue "E22" "  +0.10%  "
Set-TextValue "E23" "  +0.04%  "
Set-TextValue "D24" "81.93"
Set-TextValue "E24" "  +1.95%  "
Set-TextValue "E25" "  -4.69%  "
Set-TextValue "D26" "10.51"
Set-TextValue "E26" "  +1.89%  "
Set-TextValue "D27" "12.00"
Set-TextValue "E27" "  -2.90%  "
Set-TextValue "E28" "  +0.09%  "
Set-TextValue "E29" "  +1.76%  "
Set-TextValue "E30" "  -0.15%  "
Set-TextValue "D31" "7.20"
Set-TextValue "E31" "  +0.58%  "
Set-TextValue "E32" "  -2.34%  "
Set-TextValue "D33" "27.45"
Set-TextValue "E33" "  -0.48%  "
Set-TextValue "D34" "0.109"
Set-TextValue "E34" "  +1.94%  "
Set-TextValue "D35" "0.0₃0833"
Set-TextValue "E35" "  +4.64%  "
Set-TextValue "E36" "  -1.58%  "
Set-TextValue "D37" "5.78"
Set-TextValue "E37" "  +0.62%  "
Set-TextValue "D38" "9.29"
Set-TextValue "E38" "  +3.66%  "
Set-TextValue "D39" "50.36"
Set-TextValue "E39" "  +0.27%  "
Set-TextValue "E40" "  -3.23%  "
Set-TextValue "E41" "  +9.57%  "
Set-TextValue "E42" "  -1.74%  "
Set-TextValue "D43" "396.27"
Set-TextValue "E43" "  -6.18%  "
Set-TextValue "E44" "  -0.03%  "
Set-TextValue "D45" "39.15"
Set-TextValue "E45" "  +4.38%  "
Set-TextValue "E46" "  -3.11%  "
Set-TextValue "D47" "2.716.75"
Set-TextValue "E47" "  -2.70%  "
Set-TextValue "D48" "133.23"
Set-TextValue "E48" "  +3.61%  "
Set-TextValue "D50" "0.107"
Set-TextValue "E50" "  -0.35%  "
Set-TextValue "E51" "  +1.87%  "
